$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# PossibleVar column (D) flips from blank to YES for several rows
$ws.Range("D9").Value = "YES"   # Board
$ws.Range("D10").Value = "YES"  # Brick
$ws.Range("D25").Value = "YES"  # DoubleExterior
$ws.Range("D71").Value = "YES"  # New
$ws.Range("D74").Value = "YES"  # Other
$ws.Range("D90").Value = "YES"  # Shingles
$ws.Range("D91").Value = "YES"  # Siding
$ws.Range("D92").Value = "YES"  # Story
$ws.Range("D94").Value = "YES"  # Stucco

# PoolArea row: PossibleVar flips YES -> NO, with an explanatory note
$ws.Range("D78").Value = "NO"
$ws.Range("E78").Value = "PoolQC will be used instead"

# PoolQC row: now Examined (X) and PossibleVar = YES
$ws.Range("B79").Value = "X"
$ws.Range("D79").Value = "YES"

# PUD row: PossibleVar set to NO
$ws.Range("D81").Value = "NO"

# Restore the selection seen in the edited workbook
$ws.Activate()
$ws.Range("E78").Select()
